$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 240094.81
$ws.Cells.Item(132, 9).Value = 270420.75
$ws.Cells.Item(132, 10).Value = 45142.285
$ws.Cells.Item(132, 11).Value = 811262.25
$ws.Cells.Item(132, 12).Value = 135426.855
$ws.Cells.Item(132, 13).Value = -808732.25
$ws.Cells.Item(132, 14).Value = -140486.855

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 7938992.5
$ws.Cells.Item(138, 9).Value = 2582.7727
$ws.Cells.Item(138, 11).Value = 7748.3181
$ws.Cells.Item(138, 13).Value = -2608.3181

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 9403.571
$ws.Cells.Item(2, 9).Value = 12688
$ws.Cells.Item(2, 10).Value = 1192.5
$ws.Cells.Item(2, 11).Value = 12688
$ws.Cells.Item(2, 12).Value = 1192.5
$ws.Cells.Item(2, 13).Value = -12575
$ws.Cells.Item(2, 14).Value = -1418.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2182.5945
$ws.Cells.Item(32, 9).Value = 2118.889
$ws.Cells.Item(32, 10).Value = 2354.6
$ws.Cells.Item(32, 11).Value = 2118.889
$ws.Cells.Item(32, 12).Value = 2354.6
$ws.Cells.Item(32, 13).Value = -1831.889
$ws.Cells.Item(32, 14).Value = -2928.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1810.3438
$ws.Cells.Item(61, 9).Value = 1404.0834
$ws.Cells.Item(61, 11).Value = 1404.0834
$ws.Cells.Item(61, 13).Value = -1192.0834

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 9301.934
$ws.Cells.Item(63, 9).Value = 9823.5
$ws.Cells.Item(63, 10).Value = 2000
$ws.Cells.Item(63, 11).Value = 9823.5
$ws.Cells.Item(63, 12).Value = 2000
$ws.Cells.Item(63, 13).Value = -9137.5
$ws.Cells.Item(63, 14).Value = -3372

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 9301.934
$ws.Cells.Item(66, 9).Value = 9823.5
$ws.Cells.Item(66, 10).Value = 2000
$ws.Cells.Item(66, 11).Value = 49117.5
$ws.Cells.Item(66, 12).Value = 10000
$ws.Cells.Item(66, 13).Value = -45685.5
$ws.Cells.Item(66, 14).Value = -16864

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 10028.571
$ws.Cells.Item(88, 9).Value = 20000
$ws.Cells.Item(88, 10).Value = 6040
$ws.Cells.Item(88, 11).Value = 20000
$ws.Cells.Item(88, 12).Value = 6040
$ws.Cells.Item(88, 13).Value = -19594
$ws.Cells.Item(88, 14).Value = -6852

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 10028.571
$ws.Cells.Item(91, 9).Value = 20000
$ws.Cells.Item(91, 10).Value = 6040
$ws.Cells.Item(91, 11).Value = 20000
$ws.Cells.Item(91, 12).Value = 6040
$ws.Cells.Item(91, 13).Value = -18596
$ws.Cells.Item(91, 14).Value = -8848

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 9403.571
$ws.Cells.Item(116, 9).Value = 12688
$ws.Cells.Item(116, 10).Value = 1192.5
$ws.Cells.Item(116, 11).Value = 12688
$ws.Cells.Item(116, 12).Value = 1192.5
$ws.Cells.Item(116, 13).Value = -10394
$ws.Cells.Item(116, 14).Value = -5780.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1810.3438
$ws.Cells.Item(136, 9).Value = 1404.0834
$ws.Cells.Item(136, 11).Value = 4212.2502
$ws.Cells.Item(136, 13).Value = -1662.2502

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 9403.571
$ws.Cells.Item(3, 9).Value = 12688
$ws.Cells.Item(3, 10).Value = 1192.5
$ws.Cells.Item(3, 11).Value = 12688
$ws.Cells.Item(3, 12).Value = 1192.5
$ws.Cells.Item(3, 13).Value = -12574
$ws.Cells.Item(3, 14).Value = -1420.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 1662.1666
$ws.Cells.Item(54, 9).Value = 1662.1666
$ws.Cells.Item(54, 11).Value = 1662.1666
$ws.Cells.Item(54, 13).Value = -1178.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(57, 8).Value = 46369.152
$ws.Cells.Item(57, 9).Value = 35709
$ws.Cells.Item(57, 10).Value = 105000
$ws.Cells.Item(57, 11).Value = 35709
$ws.Cells.Item(57, 12).Value = 105000
$ws.Cells.Item(57, 13).Value = -34989
$ws.Cells.Item(57, 14).Value = -106440

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(136, 8).Value = 46369.152
$ws.Cells.Item(136, 9).Value = 35709
$ws.Cells.Item(136, 10).Value = 105000
$ws.Cells.Item(136, 11).Value = 35709
$ws.Cells.Item(136, 12).Value = 105000
$ws.Cells.Item(136, 13).Value = -30609
$ws.Cells.Item(136, 14).Value = -115200

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 1005
$ws.Cells.Item(19, 9).Value = 1005
$ws.Cells.Item(19, 11).Value = 1005
$ws.Cells.Item(19, 13).Value = -835

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(24, 8).Value = 1005
$ws.Cells.Item(24, 9).Value = 1005
$ws.Cells.Item(24, 11).Value = 1005
$ws.Cells.Item(24, 13).Value = -835

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1587.1562
$ws.Cells.Item(58, 9).Value = 926.35297
$ws.Cells.Item(58, 10).Value = 2336.0667
$ws.Cells.Item(58, 11).Value = 926.35297
$ws.Cells.Item(58, 12).Value = 2336.0667
$ws.Cells.Item(58, 13).Value = -723.35297
$ws.Cells.Item(58, 14).Value = -2742.0667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(133, 8).Value = 12747.5
$ws.Cells.Item(133, 10).Value = 12747.5
$ws.Cells.Item(133, 12).Value = 12747.5
$ws.Cells.Item(133, 14).Value = -17807.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 1587.1562
$ws.Cells.Item(136, 9).Value = 926.35297
$ws.Cells.Item(136, 10).Value = 2336.0667
$ws.Cells.Item(136, 11).Value = 2779.05891
$ws.Cells.Item(136, 12).Value = 7008.2001
$ws.Cells.Item(136, 13).Value = -229.0589100000002
$ws.Cells.Item(136, 14).Value = -12108.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2601.1758
$ws.Cells.Item(68, 9).Value = 3400.7827
$ws.Cells.Item(68, 10).Value = 1287.5358
$ws.Cells.Item(68, 11).Value = 10202.3481
$ws.Cells.Item(68, 12).Value = 3862.6074
$ws.Cells.Item(68, 13).Value = -9391.348100000001
$ws.Cells.Item(68, 14).Value = -5484.607400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 2601.1758
$ws.Cells.Item(71, 9).Value = 3400.7827
$ws.Cells.Item(71, 10).Value = 1287.5358
$ws.Cells.Item(71, 11).Value = 30607.0443
$ws.Cells.Item(71, 12).Value = 11587.8222
$ws.Cells.Item(71, 13).Value = -26551.0443
$ws.Cells.Item(71, 14).Value = -19699.8222

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 751759.5
$ws.Cells.Item(103, 10).Value = 1502519
$ws.Cells.Item(103, 12).Value = 4507557
$ws.Cells.Item(103, 14).Value = -4509315

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 2002.2858
$ws.Cells.Item(129, 9).Value = 10000
$ws.Cells.Item(129, 10).Value = 1387.0769
$ws.Cells.Item(129, 11).Value = 30000
$ws.Cells.Item(129, 12).Value = 4161.2307
$ws.Cells.Item(129, 13).Value = -25000
$ws.Cells.Item(129, 14).Value = -14161.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2309.978
$ws.Cells.Item(131, 10).Value = 2539.9507
$ws.Cells.Item(131, 12).Value = 7619.8521
$ws.Cells.Item(131, 14).Value = -17699.8521

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(54, 8).Value = 16750
$ws.Cells.Item(54, 10).Value = 5666.6665
$ws.Cells.Item(54, 12).Value = 5666.6665
$ws.Cells.Item(54, 14).Value = -6446.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 13000
$ws.Cells.Item(135, 10).Value = 13000
$ws.Cells.Item(135, 12).Value = 13000
$ws.Cells.Item(135, 14).Value = -23140

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3266.4167
$ws.Cells.Item(7, 9).Value = 2478.8
$ws.Cells.Item(7, 10).Value = 3473.6843
$ws.Cells.Item(7, 11).Value = 2478.8
$ws.Cells.Item(7, 12).Value = 3473.6843
$ws.Cells.Item(7, 13).Value = -2366.8
$ws.Cells.Item(7, 14).Value = -3697.6843

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3266.4167
$ws.Cells.Item(126, 9).Value = 2478.8
$ws.Cells.Item(126, 10).Value = 3473.6843
$ws.Cells.Item(126, 11).Value = 7436.400000000001
$ws.Cells.Item(126, 12).Value = 10421.0529
$ws.Cells.Item(126, 13).Value = -4966.400000000001
$ws.Cells.Item(126, 14).Value = -15361.0529

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

